# Updated cryptos list on Fri Jan 26 13:24:33 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Some "price" values look like plain decimal numbers (e.g. "302.48").
    # Excel's Range.Value setter auto-converts those into numeric cells,
    # which mangles trailing zeros / exact text. Force text interpretation,
    # then restore the default "Normal" style so no stray number format
    # sticks to the cell (matches the source file, which has no style
    # override on these data cells).
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "41.210.68"
$ws.Range("E2").Value = "  +3.39%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.248.60"
$ws.Range("E3").Value = "  +1.93%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - BNB
Set-TextValue "D5" "302.48"
$ws.Range("E5").Value = "  +3.02%  "

# Row 6 - Solana
Set-TextValue "D6" "90.84"
$ws.Range("E6").Value = "  +4.61%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.08%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +2.02%  "

# Row 10 - OKB
Set-TextValue "D10" "53.66"
$ws.Range("E10").Value = "  +9.26%  "

# Row 11 - Avalanche
Set-TextValue "D11" "31.81"
$ws.Range("E11").Value = "  +6.99%  "

# Row 12 - Dogecoin
Set-TextValue "D12" "0.0793"
$ws.Range("E12").Value = "  +2.60%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +3.38%  "

# Row 14 - Polkadot
Set-TextValue "D14" "6.56"
$ws.Range("E14").Value = "  +1.47%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "2.597.39"
$ws.Range("E15").Value = "  +1.90%  "

# Row 16 - now WrappedEther (was Chainlink)
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D16" "2.403.03"
$ws.Range("E16").Value = "  +8.65%  "

# Row 17 - now Chainlink (was WrappedEther)
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D17" "14.10"
$ws.Range("E17").Value = "  +2.85%  "

# Row 18 - Polygon
Set-TextValue "D18" "0.749"
$ws.Range("E18").Value = "  +3.14%  "

# Row 19 - WrappedBTC
Set-TextValue "D19" "41.139.28"
$ws.Range("E19").Value = "  +3.39%  "

# Row 20 - InternetComputer(DFINITY)
Set-TextValue "D20" "11.89"
$ws.Range("E20").Value = "  +5.68%  "

# Row 21 - ShibaInu
Set-TextValue "D21" "0.0₃0902"
$ws.Range("E21").Value = "  +2.26%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +1.57%  "

# Row 23 - Litecoin
Set-TextValue "D23" "66.85"
$ws.Range("E23").Value = "  +2.67%  "

# Row 24 - BitcoinCash
Set-TextValue "D24" "240.30"
$ws.Range("E24").Value = "  +1.80%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +4.29%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  -0.03%  "

# Row 27 - ImmutableX
$ws.Range("E27").Value = "  +2.59%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "23.77"
$ws.Range("E28").Value = "  +5.93%  "

# Row 29 - Cosmos
Set-TextValue "D29" "9.56"
$ws.Range("E29").Value = "  +4.43%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  -2.53%  "

# Row 31 - Monero
Set-TextValue "D31" "158.90"
$ws.Range("E31").Value = "  +2.38%  "

# Row 32 - InjectiveProtocol
Set-TextValue "D32" "33.22"

# Row 33 - FirstDigitalUSD
$ws.Range("E33").Value = "  +0.01%  "

# Row 34 - Filecoin
Set-TextValue "D34" "5.17"
$ws.Range("E34").Value = "  +5.83%  "

# Row 35 - Hedera
$ws.Range("E35").Value = "  +2.92%  "

# Row 36 - LidoDAOToken
$ws.Range("E36").Value = "  +7.27%  "

# Row 37 - WEMIXToken
$ws.Range("E37").Value = "  +0.62%  "

# Row 38 - now Celestia (was Kaspa)
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D38" "16.48"
$ws.Range("E38").Value = "  +6.69%  "

# Row 39 - now Kaspa (was Stellar)
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D39" "0.104"
$ws.Range("E39").Value = "  +6.26%  "

# Row 40 - now Stellar (was Celestia)
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D40" "0.115"
$ws.Range("E40").Value = "  +3.25%  "

# Row 41 - ARBITRUM
$ws.Range("E41").Value = "  +7.20%  "

# Row 42 - RenderToken
Set-TextValue "D42" "3.91"
$ws.Range("E42").Value = "  +4.77%  "

# Row 43 - Maker
Set-TextValue "D43" "2.072.66"
$ws.Range("E43").Value = "  -2.19%  "

# Row 44 - EnergySwap
Set-TextValue "D44" "20.13"
$ws.Range("E44").Value = "  +13.59%  "

# Row 45 - VeChain
$ws.Range("E45").Value = "  +3.88%  "

# Row 46 - FraxShare
Set-TextValue "D46" "10.24"
$ws.Range("E46").Value = "  +6.53%  "

# Row 47 - NEARProtocol
Set-TextValue "D47" "2.94"
$ws.Range("E47").Value = "  +10.90%  "

# Row 48 - ApeXProtocol
Set-TextValue "D48" "1.99"
$ws.Range("E48").Value = "  -5.60%  "

# Row 49 - RocketPoolETH
Set-TextValue "D49" "2.469.02"
$ws.Range("E49").Value = "  +2.13%  "

# Row 50 - Stacks
Set-TextValue "D50" "1.50"
$ws.Range("E50").Value = "  +3.73%  "

# Row 51 - TrustWalletToken
$ws.Range("E51").Value = "  +4.33%  "
